# Insert a new row at position 220 (this shifts existing rows 220..312 down to 221..313)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new data record.
$ws.Cells.Item(220, 1).Value = 9
$ws.Cells.Item(220, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(220, 3).Value = "Metropolitana"
$ws.Cells.Item(220, 4).Value = 44837
$ws.Cells.Item(220, 5).Value = 13
$ws.Cells.Item(220, 6).Value = 100112001
$ws.Cells.Item(220, 7).Value = "Berenjena"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 70
$ws.Cells.Item(220, 11).Value = 12000
$ws.Cells.Item(220, 12).Value = 13000
$ws.Cells.Item(220, 13).Value = 12500
$ws.Cells.Item(220, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(220, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(220, 16).Value = 250
$ws.Cells.Item(220, 17).Value = 50
$ws.Cells.Item(220, 18).Value = "Hortaliza"
